$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Trae Young", "PG", "Atlanta Hawks")
    ,@("LeBron James", "SF,PF", "Los Angeles Lakers")
    ,@("Matas Buzelis", "SF,PF", "Chicago Bulls")
    ,@("Jimmy Butler III", "SF,PF", "Golden State Warriors")
    ,@("Myles Turner", "C", "Indiana Pacers")
    ,@("Nicolas Claxton", "C", "Brooklyn Nets")
    ,@("Jonas Valanciunas", "C", "Sacramento Kings")
    ,@("Desmond Bane", "SG,SF", "Memphis Grizzlies")
    ,@("Santi Aldama", "PF,C", "Memphis Grizzlies")
    ,@("Walker Kessler", "C", "Utah Jazz")
    ,@("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
    ,@("Miles McBride", "PG,SG", "New York Knicks")
    ,@("Devin Booker", "PG,SG", "Phoenix Suns")
    ,@("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder")
    ,@("Brandon Ingram", "SG,SF,PF", "Toronto Raptors")
    ,@("Norman Powell", "SG,SF", "LA Clippers")
    ,@("Jalen Brunson", "PG", "New York Knicks")
    ,@("Immanuel Quickley", "PG,SG", "Toronto Raptors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
